$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9: a new "定制" (Customize) button/card is added - text updates.
# Order matters: set C9 first so its shared string is inserted before B9's,
# matching the shared-string table ordering of the target workbook.
$ws.Range("C9").Value = "查看我的副本定制卡组(D)"
$ws.Range("B9").Value = "定制"

# ShowInDungeon (column H) flips from true to false for rows 10 and 17-22.
# Assigning the literal string "false" directly would make Excel coerce the
# cell to a native boolean type, so instead copy an existing text "false"
# cell (G9) and paste-special its value, which preserves the shared-string
# text type and keeps the original cell style.
foreach ($r in 10, 17, 18, 19, 20, 21, 22) {
    $ws.Range("G9").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# Update the selected range shown in the sheet view.
$ws.Range("H17:H22").Select() | Out-Null
$excel.CutCopyMode = 0
